# Update the "department" column (C) values on the "courses" sheet.
# Previously every row used the single literal school name
# "BRANSON SCHOOL OF BUSINESS AND TECHNOLOGY"; now each course row is
# categorised into its proper department/category.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C4").Value   = "Accounting"
$ws.Range("C5:C7").Value   = "Information Technology"
$ws.Range("C8:C9").Value   = "Logistics"
$ws.Range("C10:C14").Value = "Management"
$ws.Range("C15:C16").Value = "Graduate Studies"
$ws.Range("C17:C22").Value = "Packages"
